$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a Number by Excel
# (single-decimal-point numeric strings) are forced back to Text first, so they
# round-trip exactly like the original inline-string cells.

$ws.Range("D2").Value = "26.026.92"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "1.650.21"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.32"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5193"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2634"
$ws.Range("E8").Value = "  +0.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06321"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07648"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.579"
$ws.Range("E12").Value = "  +2.43%  "

$ws.Range("D13").Value = "1.652.64"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").Value = "1.877.43"
$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5593"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "0.0₅8139"
$ws.Range("E16").Value = "  +2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.20"

$ws.Range("D18").Value = "26.033.03"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.616"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("E21").Value = "  +4.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.57"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.901"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.48"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1185"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.190"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.86"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05381"
$ws.Range("E30").Value = "  -4.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.269"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.452"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.351"
$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.553"
$ws.Range("E34").Value = "  -2.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.418"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9459"
$ws.Range("E37").Value = "  +1.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5631"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01577"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.869"
$ws.Range("E40").Value = "  -0.95%  "

$ws.Range("D42").Value = "1.029.51"
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8260"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.74"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").Value = "1.787.16"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("E46").Value = "  +5.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.37"
$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9989"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.912"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05141"
$ws.Range("E51").Value = "  -3.78%  "
